# IST price update 2025-12-20 22:45
# A new price-snapshot column is inserted as the new column B (pushing the
# previous B..H snapshot columns to C..I). The header row gets a new
# timestamp in B1; most data rows shift their old B..H values right by one
# column and leave the new B cell blank, but a handful of rows in the
# source sheet were captured slightly differently and only gained a extra
# I-column echo of their H value (no shift, no blank). Row 19's H/I pair
# also carries a genuine price change (H overwritten, old value preserved
# in I). We reproduce the exact resulting grid explicitly, cell by cell,
# rather than relying on a generic column-insert (whose behaviour does not
# match every row here).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: give it the same width as the other price columns ---
$ws.Columns("I").ColumnWidth = 20.17

# --- Row 1 (header/timestamps): new B1, then old B1..H1 shift to C1..I1 ---
$ws.Range("I1").Value = $ws.Range("H1").Value()
$ws.Range("H1").Value = $ws.Range("G1").Value()
$ws.Range("G1").Value = $ws.Range("F1").Value()
$ws.Range("F1").Value = $ws.Range("E1").Value()
$ws.Range("E1").Value = $ws.Range("D1").Value()
$ws.Range("D1").Value = $ws.Range("C1").Value()
$ws.Range("C1").Value = $ws.Range("B1").Value()
$ws.Range("B1").Value = "2025-12-21 04:09"

# Give the newly-created I1 header the same formatting as the rest of row 1
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rows that shift B..H -> C..I and leave the new B blank ---
$shiftRows = @(3,4,5,6,7,8,9,10,11,12,13,15,16,17,18,20,22,23,24,25,26)
foreach ($r in $shiftRows) {
    $ws.Cells.Item($r, 9).Value = $ws.Cells.Item($r, 8).Value()   # I = old H
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($r, 7).Value()   # H = old G
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 6).Value()   # G = old F
    $ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 5).Value()   # F = old E
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 4).Value()   # E = old D
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 3).Value()   # D = old C
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 2).Value()   # C = old B
    $ws.Cells.Item($r, 2).ClearContents()                          # B = blank
}

# --- Rows 2, 14, 21: no shift, just echo the H price into the new I cell ---
foreach ($r in @(2,14,21)) {
    $ws.Cells.Item($r, 9).Value = $ws.Cells.Item($r, 8).Value()
}

# --- Row 19: price change captured between H and I (no shift elsewhere) ---
$ws.Range("I19").Value = $ws.Range("H19").Value()   # preserve old H (2997) in I
$ws.Range("H19").Value = $ws.Range("G19").Value()   # H takes G's price (1497)
